$d = $word.ActiveDocument

# 1. Update the caption paragraph describing the table.
$caption = $d.Paragraphs.Item(2)
$caption.Range.Text = "Effects of salinity and nutrients on community composition (permutations restricted within site)."

# 2. Update the PERMANOVA results table values.
$t = $d.Tables.Item(1)

# Row 2 = "Model"
$t.Cell(2, 2).Range.Text = "5"
$t.Cell(2, 3).Range.Text = "4.792754"
$t.Cell(2, 4).Range.Text = "0.1820246"
$t.Cell(2, 5).Range.Text = "3.694009"
$t.Cell(2, 6).Range.Text = "0.005"

# Row 3 = "Residual"
$t.Cell(3, 2).Range.Text = "83"
$t.Cell(3, 3).Range.Text = "21.537497"
$t.Cell(3, 4).Range.Text = "0.8179754"

# Row 4 = "Total"
$t.Cell(4, 3).Range.Text = "26.330251"
